$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.829.23'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.084.76'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.38'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.08'
$ws.Range("E7").Value = '  +3.06%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0789'
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.106'
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.392.31'
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.29'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.770'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.29'
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.082.80'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.744.56'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.38'
$ws.Range("E20").Value = '  +2.58%  '
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.02'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.37'
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.01'
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.53'
$ws.Range("E30").Value = '  +2.15%  '
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("E32").Value = '  +3.89%  '
$ws.Range("E33").Value = '  +4.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0628'
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("E36").Value = '  +2.08%  '
$ws.Range("E37").Value = '  +2.46%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0993'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.55'
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.45'
$ws.Range("E43").Value = '  +7.49%  '
$ws.Range("E44").Value = '  +1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.460.01'
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.12'
$ws.Range("E47").Value = '  +6.39%  '
$ws.Range("E48").Value = '  +4.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.41'
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("E50").Value = '  +2.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.277.06'
$ws.Range("E51").Value = '  +0.71%  '
